# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy the formatting used by the other header cells (bold,
# bordered, centered style already in the workbook) and set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data cell H2: numeric value for the single data row.
$ws.Range("H2").Value = 1
